$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 13-40 (values shifted due to inserted/reordered weekly records) ---
$ws.Cells.Item(13, 4).Value = 44775
$ws.Cells.Item(13, 10).Value = 60
$ws.Cells.Item(13, 11).Value = 14000
$ws.Cells.Item(13, 12).Value = 15000
$ws.Cells.Item(13, 13).Value = 14500
$ws.Cells.Item(13, 16).Value = 362
$ws.Cells.Item(14, 4).Value = 44467
$ws.Cells.Item(14, 10).Value = 160
$ws.Cells.Item(14, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(15, 4).Value = 44510
$ws.Cells.Item(16, 4).Value = 44516
$ws.Cells.Item(16, 10).Value = 120
$ws.Cells.Item(16, 11).Value = 11000
$ws.Cells.Item(16, 12).Value = 12000
$ws.Cells.Item(16, 13).Value = 11500
$ws.Cells.Item(16, 16).Value = 288
$ws.Cells.Item(17, 4).Value = 44468
$ws.Cells.Item(17, 10).Value = 60
$ws.Cells.Item(17, 11).Value = 12000
$ws.Cells.Item(17, 12).Value = 13000
$ws.Cells.Item(17, 13).Value = 12500
$ws.Cells.Item(17, 16).Value = 312
$ws.Cells.Item(18, 4).Value = 44491
$ws.Cells.Item(18, 10).Value = 100
$ws.Cells.Item(19, 4).Value = 44505
$ws.Cells.Item(19, 10).Value = 120
$ws.Cells.Item(19, 11).Value = 11000
$ws.Cells.Item(19, 12).Value = 12000
$ws.Cells.Item(19, 13).Value = 11500
$ws.Cells.Item(19, 16).Value = 288
$ws.Cells.Item(20, 4).Value = 44455
$ws.Cells.Item(20, 10).Value = 100
$ws.Cells.Item(20, 11).Value = 13000
$ws.Cells.Item(20, 12).Value = 14000
$ws.Cells.Item(20, 13).Value = 13500
$ws.Cells.Item(20, 16).Value = 338
$ws.Cells.Item(21, 4).Value = 44435
$ws.Cells.Item(21, 10).Value = 120
$ws.Cells.Item(21, 11).Value = 14000
$ws.Cells.Item(21, 12).Value = 15000
$ws.Cells.Item(21, 13).Value = 14500
$ws.Cells.Item(21, 16).Value = 362
$ws.Cells.Item(22, 4).Value = 44498
$ws.Cells.Item(22, 10).Value = 60
$ws.Cells.Item(22, 11).Value = 10500
$ws.Cells.Item(22, 12).Value = 11000
$ws.Cells.Item(22, 13).Value = 10750
$ws.Cells.Item(22, 16).Value = 269
$ws.Cells.Item(23, 4).Value = 44432
$ws.Cells.Item(23, 8).Value = 'Madrigal'
$ws.Cells.Item(23, 10).Value = 120
$ws.Cells.Item(23, 11).Value = 14000
$ws.Cells.Item(23, 12).Value = 15000
$ws.Cells.Item(23, 13).Value = 14500
$ws.Cells.Item(23, 14).Value = '$/caja 40 unidades'
$ws.Cells.Item(23, 15).Value = 'Provincia del Elquí'
$ws.Cells.Item(23, 16).Value = 362
$ws.Cells.Item(23, 17).Value = 40
$ws.Cells.Item(24, 4).Value = 44762
$ws.Cells.Item(24, 8).Value = 'Argentina(o)'
$ws.Cells.Item(24, 10).Value = 60
$ws.Cells.Item(24, 11).Value = 15000
$ws.Cells.Item(24, 12).Value = 16000
$ws.Cells.Item(24, 13).Value = 15500
$ws.Cells.Item(24, 14).Value = '$/caja 50 unidades'
$ws.Cells.Item(24, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(24, 16).Value = 310
$ws.Cells.Item(24, 17).Value = 50
$ws.Cells.Item(25, 4).Value = 44454
$ws.Cells.Item(25, 8).Value = 'Madrigal'
$ws.Cells.Item(25, 10).Value = 120
$ws.Cells.Item(25, 11).Value = 13000
$ws.Cells.Item(25, 12).Value = 14000
$ws.Cells.Item(25, 13).Value = 13500
$ws.Cells.Item(25, 14).Value = '$/caja 40 unidades'
$ws.Cells.Item(25, 15).Value = 'Provincia del Elquí'
$ws.Cells.Item(25, 16).Value = 338
$ws.Cells.Item(25, 17).Value = 40
$ws.Cells.Item(26, 4).Value = 44753
$ws.Cells.Item(26, 8).Value = 'Argentina(o)'
$ws.Cells.Item(26, 10).Value = 60
$ws.Cells.Item(26, 11).Value = 16000
$ws.Cells.Item(26, 12).Value = 17000
$ws.Cells.Item(26, 13).Value = 16500
$ws.Cells.Item(26, 14).Value = '$/caja 50 unidades'
$ws.Cells.Item(26, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(26, 16).Value = 330
$ws.Cells.Item(26, 17).Value = 50
$ws.Cells.Item(27, 4).Value = 44420
$ws.Cells.Item(27, 10).Value = 120
$ws.Cells.Item(27, 11).Value = 13000
$ws.Cells.Item(27, 12).Value = 14000
$ws.Cells.Item(27, 13).Value = 13500
$ws.Cells.Item(27, 16).Value = 338
$ws.Cells.Item(28, 4).Value = 44503
$ws.Cells.Item(28, 10).Value = 160
$ws.Cells.Item(28, 11).Value = 11000
$ws.Cells.Item(28, 12).Value = 12000
$ws.Cells.Item(28, 13).Value = 11500
$ws.Cells.Item(28, 16).Value = 288
$ws.Cells.Item(29, 4).Value = 44417
$ws.Cells.Item(29, 11).Value = 15000
$ws.Cells.Item(29, 12).Value = 16000
$ws.Cells.Item(29, 13).Value = 15500
$ws.Cells.Item(29, 16).Value = 388
$ws.Cells.Item(30, 4).Value = 44515
$ws.Cells.Item(30, 11).Value = 11000
$ws.Cells.Item(30, 12).Value = 12000
$ws.Cells.Item(30, 13).Value = 11500
$ws.Cells.Item(30, 16).Value = 288
$ws.Cells.Item(31, 4).Value = 44427
$ws.Cells.Item(31, 11).Value = 13000
$ws.Cells.Item(31, 12).Value = 14000
$ws.Cells.Item(31, 13).Value = 13500
$ws.Cells.Item(31, 16).Value = 338
$ws.Cells.Item(32, 4).Value = 44494
$ws.Cells.Item(32, 10).Value = 120
$ws.Cells.Item(33, 4).Value = 44487
$ws.Cells.Item(33, 10).Value = 100
$ws.Cells.Item(34, 4).Value = 44484
$ws.Cells.Item(34, 10).Value = 120
$ws.Cells.Item(35, 4).Value = 44488
$ws.Cells.Item(35, 10).Value = 100
$ws.Cells.Item(36, 4).Value = 44496
$ws.Cells.Item(36, 8).Value = 'Madrigal'
$ws.Cells.Item(36, 9).Value = 'Primera'
$ws.Cells.Item(36, 11).Value = 11000
$ws.Cells.Item(36, 12).Value = 12000
$ws.Cells.Item(36, 13).Value = 11500
$ws.Cells.Item(36, 16).Value = 288
$ws.Cells.Item(37, 4).Value = 44399
$ws.Cells.Item(37, 8).Value = 'Española'
$ws.Cells.Item(37, 9).Value = 'Segunda'
$ws.Cells.Item(37, 11).Value = 15500
$ws.Cells.Item(37, 12).Value = 16000
$ws.Cells.Item(37, 13).Value = 15750
$ws.Cells.Item(37, 15).Value = 'Provincia del Elquí'
$ws.Cells.Item(37, 16).Value = 394
$ws.Cells.Item(38, 4).Value = 44425
$ws.Cells.Item(38, 11).Value = 14000
$ws.Cells.Item(38, 12).Value = 15000
$ws.Cells.Item(38, 13).Value = 14500
$ws.Cells.Item(38, 15).Value = 'Región del Maule'
$ws.Cells.Item(38, 16).Value = 362
$ws.Cells.Item(39, 4).Value = 44512
$ws.Cells.Item(39, 10).Value = 120
$ws.Cells.Item(39, 11).Value = 11000
$ws.Cells.Item(39, 12).Value = 12000
$ws.Cells.Item(39, 13).Value = 11500
$ws.Cells.Item(39, 16).Value = 288
$ws.Cells.Item(40, 4).Value = 44453
$ws.Cells.Item(40, 10).Value = 160
$ws.Cells.Item(40, 11).Value = 12500
$ws.Cells.Item(40, 12).Value = 13000
$ws.Cells.Item(40, 13).Value = 12750
$ws.Cells.Item(40, 16).Value = 319

# --- Insert new row 41 (new weekly record) ---
$ws.Cells.Item(41, 1).Value = 7
$ws.Cells.Item(41, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(41, 3).Value = 'Ñuble'
$ws.Cells.Item(41, 4).Value = 44489
$ws.Cells.Item(41, 5).Value = 16
$ws.Cells.Item(41, 6).Value = 100112013
$ws.Cells.Item(41, 7).Value = 'Alcachofa'
$ws.Cells.Item(41, 8).Value = 'Madrigal'
$ws.Cells.Item(41, 9).Value = 'Primera'
$ws.Cells.Item(41, 10).Value = 120
$ws.Cells.Item(41, 11).Value = 11000
$ws.Cells.Item(41, 12).Value = 12000
$ws.Cells.Item(41, 13).Value = 11500
$ws.Cells.Item(41, 14).Value = '$/caja 40 unidades'
$ws.Cells.Item(41, 15).Value = 'Provincia del Elquí'
$ws.Cells.Item(41, 16).Value = 288
$ws.Cells.Item(41, 17).Value = 40
$ws.Cells.Item(41, 18).Value = 'Hortaliza'

$ws.Cells.Item(41, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

